$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OverallRebateEfficiency")
$ws2 = $wb.Worksheets.Item("PSA_LOLO")

# Append new weekly rows 25-30 to the OverallRebateEfficiency sheet
$weeks = @("Week_25", "Week_26", "Week_27", "Week_28", "Week_29", "Week_30")
$values = @(0.7175, 0.7577, 0.7187, 0.7243, 0.7219, 0.7721)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $r = 25 + $i
    $ws1.Cells.Item($r, 1).Value = $weeks[$i]
    $ws1.Cells.Item($r, 2).Value = $values[$i]
}

# Update the view so the newly entered area is visible/selected
$ws1.Range("A31").Select() | Out-Null

# Update PSA_LOLO sheet values
$ws2.Range("A2").Value = 40297
$ws2.Range("B2").Value = 22479

Write-Host "edits applied"
